$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Assets")

$newItems = @(
    "PwUpdated_BNSTESTERB",
    "PwUpdated_BNSTESTERC",
    "TestData_ProductID",
    "TestData_ProductType",
    "TestData_StateCode"
)

$row = 4
foreach ($item in $newItems) {
    $ws.Range("A$row").Value = $item
    $ws.Range("B$row").Value = $item
    $row++
}
